$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" '28.488.52'
Set-TextValue "E2" '  +1.93%  '
Set-TextValue "D3" '1.909.75'
Set-TextValue "E3" '  +5.34%  '
Set-TextValue "D4" '1.001'
Set-TextValue "E4" '  -0.18%  '
Set-TextValue "D5" '313.84'
Set-TextValue "E5" '  +1.16%  '
Set-TextValue "D6" '1.001'
Set-TextValue "D7" '0.5051'
Set-TextValue "E7" '  +1.31%  '
Set-TextValue "D8" '0.3945'
Set-TextValue "E8" '  +1.44%  '
Set-TextValue "D9" '0.09799'
Set-TextValue "E9" '  +1.13%  '
Set-TextValue "D10" '1.162'
Set-TextValue "E10" '  +5.41%  '
Set-TextValue "D11" '41.55'
Set-TextValue "E11" '  +1.31%  '
Set-TextValue "D12" '6.547'
Set-TextValue "E12" '  +1.90%  '
Set-TextValue "D13" '21.12'
Set-TextValue "E13" '  +3.01%  '
Set-TextValue "D14" '1.911.51'
Set-TextValue "E14" '  +5.45%  '
Set-TextValue "D15" '7.585'
Set-TextValue "E15" '  +3.82%  '
Set-TextValue "D16" '1.001'
Set-TextValue "E16" '  -0.06%  '
Set-TextValue "D17" '0.00001144'
Set-TextValue "E17" '  +1.71%  '
Set-TextValue "D18" '93.86'
Set-TextValue "E18" '  +1.34%  '
Set-TextValue "D19" '0.06655'
Set-TextValue "E19" '  -0.09%  '
Set-TextValue "D20" '18.07'
Set-TextValue "E20" '  +5.51%  '
Set-TextValue "E21" '  -0.04%  '
Set-TextValue "D22" '6.294'
Set-TextValue "E22" '  +6.34%  '
Set-TextValue "D23" '28.547.19'
Set-TextValue "E23" '  +1.88%  '
Set-TextValue "D24" '11.46'
Set-TextValue "E24" '  +3.18%  '
Set-TextValue "D25" '2.282'
Set-TextValue "E25" '  +1.88%  '
Set-TextValue "D26" '2.735'
Set-TextValue "E26" '  +14.35%  '
Set-TextValue "D27" '2.127.23'
Set-TextValue "E27" '  +5.21%  '
Set-TextValue "E28" '  +3.79%  '
Set-TextValue "D29" '159.77'
Set-TextValue "E29" '  +0.47%  '
Set-TextValue "D30" '128.98'
Set-TextValue "E30" '  +0.71%  '
Set-TextValue "D31" '1.106'
Set-TextValue "E31" '  +6.38%  '
Set-TextValue "D32" '0.1071'
Set-TextValue "E32" '  +0.98%  '
Set-TextValue "D33" '5.707'
Set-TextValue "E33" '  +2.43%  '
Set-TextValue "D34" '3.639'
Set-TextValue "D35" '9.943'
Set-TextValue "E35" '  +10.66%  '
Set-TextValue "D36" '0.06791'
Set-TextValue "E36" '  +1.07%  '
Set-TextValue "D37" '0.02447'
Set-TextValue "E37" '  +5.01%  '
Set-TextValue "D38" '0.2239'
Set-TextValue "E38" '  +4.81%  '
Set-TextValue "D39" '5.093'
Set-TextValue "E39" '  +3.05%  '
Set-TextValue "B40" 'TheSandbox'
Set-TextValue "C40" 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue "D40" '0.6442'
Set-TextValue "E40" '  +4.07%  '
Set-TextValue "B41" 'Aptos'
Set-TextValue "C41" 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue "D41" '11.62'
Set-TextValue "E41" '  +3.32%  '
Set-TextValue "D42" '1.194'
Set-TextValue "E42" '  +4.12%  '
Set-TextValue "E43" '  +0.00%  '
Set-TextValue "D44" '13.75'
Set-TextValue "E44" '  +4.26%  '
Set-TextValue "D45" '0.6116'
Set-TextValue "E45" '  +4.16%  '
Set-TextValue "D46" '1.282'
Set-TextValue "E46" '  +0.28%  '
Set-TextValue "D47" '3.667'
Set-TextValue "E47" '  -0.60%  '
Set-TextValue "D48" '2.048'
Set-TextValue "E48" '  +5.64%  '
Set-TextValue "D49" '125.32'
Set-TextValue "E49" '  +1.87%  '
Set-TextValue "D50" '1.212'
Set-TextValue "E50" '  +2.77%  '
Set-TextValue "D51" '78.45'
Set-TextValue "E51" '  +6.38%  '
